# Dll3-Notch1.xlsx -- refresh the LR-pair table with the new TPM run.
# A 4th sending/target cluster ("Resolving-Mac") joins the existing ECs/FAPs/MuSCs
# set, so every row gets new Ligand-/Receptor-/Edge-derived figures, and the table
# grows from 12 data rows to 16 (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
     (2, "ECs", "Dll3", "Notch1", "ECs", 3, 1, 1.305319666666667, 3.915959, 0.4586709810613518, 0.4586709810613519, 3, 1, 58.95713633333333, 176.871409, 0.4863146960083892, 0.4863146960083893, 76.95790954624789, 692.621185916231, 0.223058438722721, 0.223058438722721)
    ,(3, "ECs", "Dll3", "Notch1", "FAPs", 3, 1, 1.305319666666667, 3.915959, 0.4586709810613518, 0.4586709810613519, 3, 1, 8.99153, 26.97459, 0.07416766570679004, 0.07416766570679005, 11.73682094242333, 105.63138848181, 0.03401855599276377, 0.03401855599276378)
    ,(4, "ECs", "Dll3", "Notch1", "MuSCs", 3, 1, 1.305319666666667, 3.915959, 0.4586709810613518, 0.4586709810613519, 3, 1, 42.51661933333333, 127.549858, 0.3507032073181665, 0.3507032073181665, 55.49777937598022, 499.480014383822, 0.1608573841619861, 0.1608573841619861)
    ,(5, "ECs", "Dll3", "Notch1", "Resolving-Mac", 3, 1, 1.305319666666667, 3.915959, 0.4586709810613518, 0.4586709810613519, 3, 1, 10.76719366666667, 32.301581, 0.0888144309666542, 0.08881443096665421, 14.05462964790878, 126.491666831179, 0.04073660218388099, 0.040736602183881)
    ,(6, "FAPs", "Dll3", "Notch1", "ECs", 3, 1, 0.983774, 2.951322, 0.3456843539904149, 0.3456843539904149, 3, 1, 58.95713633333333, 176.871409, 0.4863146960083892, 0.4863146960083893, 58.00049783918867, 522.004480552698, 0.168111381525705, 0.1681113815257051)
    ,(7, "FAPs", "Dll3", "Notch1", "FAPs", 3, 1, 0.983774, 2.951322, 0.3456843539904149, 0.3456843539904149, 3, 1, 8.99153, 26.97459, 0.07416766570679004, 0.07416766570679005, 8.84563343422, 79.61070090798, 0.02563860160682876, 0.02563860160682877)
    ,(8, "FAPs", "Dll3", "Notch1", "MuSCs", 3, 1, 0.983774, 2.951322, 0.3456843539904149, 0.3456843539904149, 3, 1, 42.51661933333333, 127.549858, 0.3507032073181665, 0.3507032073181665, 41.82674466803066, 376.440702012276, 0.1212326116641469, 0.1212326116641469)
    ,(9, "FAPs", "Dll3", "Notch1", "Resolving-Mac", 3, 1, 0.983774, 2.951322, 0.3456843539904149, 0.3456843539904149, 3, 1, 10.76719366666667, 32.301581, 0.0888144309666542, 0.08881443096665421, 10.59248518223133, 95.332366640082, 0.03070175919373416, 0.03070175919373416)
    ,(10, "MuSCs", "Dll3", "Notch1", "ECs", 1, 0.3333333333333333, 0.1949326666666667, 0.584798, 0.0684965987597716, 0.06849659875977161, 3, 1, 58.95713633333333, 176.871409, 0.4863146960083892, 0.4863146960083893, 11.49267180448689, 103.434046240382, 0.03331090260346693, 0.03331090260346695)
    ,(11, "MuSCs", "Dll3", "Notch1", "FAPs", 1, 0.3333333333333333, 0.1949326666666667, 0.584798, 0.0684965987597716, 0.06849659875977161, 3, 1, 8.99153, 26.97459, 0.07416766570679004, 0.07416766570679005, 1.752742920313333, 15.77468628282, 0.005080232838866869, 0.005080232838866871)
    ,(12, "MuSCs", "Dll3", "Notch1", "MuSCs", 1, 0.3333333333333333, 0.1949326666666667, 0.584798, 0.0684965987597716, 0.06849659875977161, 3, 1, 42.51661933333333, 127.549858, 0.3507032073181665, 0.3507032073181665, 8.287877984298222, 74.590901858684, 0.02402197687543744, 0.02402197687543745)
    ,(13, "MuSCs", "Dll3", "Notch1", "Resolving-Mac", 1, 0.3333333333333333, 0.1949326666666667, 0.584798, 0.0684965987597716, 0.06849659875977161, 3, 1, 10.76719366666667, 32.301581, 0.0888144309666542, 0.08881443096665421, 2.098877773959777, 18.889899965638, 0.006083486442000346, 0.006083486442000349)
    ,(14, "Resolving-Mac", "Dll3", "Notch1", "ECs", 3, 1, 0.3618473333333334, 1.085542, 0.1271480661884616, 0.1271480661884616, 3, 1, 58.95713633333333, 176.871409, 0.4863146960083892, 0.4863146960083893, 21.33348256318644, 192.001343068678, 0.06183397315649627, 0.06183397315649628)
    ,(15, "Resolving-Mac", "Dll3", "Notch1", "FAPs", 3, 1, 0.3618473333333334, 1.085542, 0.1271480661884616, 0.1271480661884616, 3, 1, 8.99153, 26.97459, 0.07416766570679004, 0.07416766570679005, 3.253561153086666, 29.28205037778, 0.009430275268330636, 0.009430275268330637)
    ,(16, "Resolving-Mac", "Dll3", "Notch1", "MuSCs", 3, 1, 0.3618473333333334, 1.085542, 0.1271480661884616, 0.1271480661884616, 3, 1, 42.51661933333333, 127.549858, 0.3507032073181665, 0.3507032073181665, 15.38452532811511, 138.460727953036, 0.04459123461659601, 0.04459123461659602)
    ,(17, "Resolving-Mac", "Dll3", "Notch1", "Resolving-Mac", 3, 1, 0.3618473333333334, 1.085542, 0.1271480661884616, 0.1271480661884616, 3, 1, 10.76719366666667, 32.301581, 0.0888144309666542, 0.08881443096665421, 3.896080315766889, 35.064722841902, 0.01129258314703871, 0.01129258314703871)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($i = 1; $i -lt $r.Length; $i++) {
        $col = [char](64 + $i)   # 1 -> A, 2 -> B, ... 20 -> T
        $ws.Range("$col$rowNum").Value = $r[$i]
    }
}
